$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell M6: "Ponderado", bold, left border, same fill flag as A9/B9 style
$ws.Range("A9").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = "Ponderado"
$ws.Range("M6").VerticalAlignment = -4107
$ws.Range("M6").WrapText = $false
$ws.Range("M6").Borders.Item(7).LineStyle = 1

# New column M formulas: weighted total = L * 10 / 100
$ws.Range("M7").Formula = "=L7*10/100"
$ws.Range("M8:M12").Formula = "=L8*10/100"
$ws.Range("M7:M12").NumberFormat = "0"

# column width for M
$ws.Columns.Item(13).ColumnWidth = 11.28515625

# Sheet view changes
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("M12").Select()
